# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1) Bump the "Date" metadata value on the Metadata sheet.
# 2) Add a new mapping column ("Mapping: Spécification métier vers
#    l'extension ROR NbPermanentSocialHelpPlace") to the Elements sheet,
#    with the business field name filled in on the Extension.value[x] row.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 : Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2) Elements sheet: new mapping column (column AL / 38) ---------------
$ws = $wb.Worksheets.Item("Elements")

$mappingHeader = "Mapping: Spécification métier vers l'extension ROR NbPermanentSocialHelpPlace"

# Header cell (row 1) - copy the header style from the neighbouring
# "Mapping: RIM Mapping" column (AK1) so the new column looks the same.
$ws.Range("AL1").Value = $mappingHeader
$ws.Range("AK1").Copy()
$ws.Range("AL1").PasteSpecial(-4122)
$ws.Range("AL1").Value = $mappingHeader

# Data cells (rows 2-6). Rows 2-5 stay blank for this mapping; row 6
# (Extension.value[x]) carries the business field name.
$ws.Range("AL2").Value = "'"
$ws.Range("AL3").Value = "'"
$ws.Range("AL4").Value = "'"
$ws.Range("AL5").Value = "'"
$ws.Range("AL6").Value = "nbPlaceAideSocialPermanent"

$ws.Range("AK2:AK5").Copy()
$ws.Range("AL2:AL5").PasteSpecial(-4122)

$ws.Range("AK6").Copy()
$ws.Range("AL6").PasteSpecial(-4122)
$ws.Range("AL6").Value = "nbPlaceAideSocialPermanent"

$ws.Columns.Item(38).ColumnWidth = 85.67578125
